# Apply updated crypto price/volume values for Sheet1 (cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.166.57"
$ws.Range("E2").Value = "  +2.58%  "
$ws.Range("D3").Value = "2.056.36"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'229.53"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "'0.614"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").Value = "'61.10"
$ws.Range("E7").Value = "  +8.89%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("E10").Value = "  +2.65%  "
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("D12").Value = "'14.78"
$ws.Range("E12").Value = "  +3.27%  "
$ws.Range("D13").Value = "2.357.46"
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("D14").Value = "'21.13"
$ws.Range("E14").Value = "  +5.11%  "
$ws.Range("D15").Value = "'5.35"
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D16").Value = "'0.758"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").Value = "2.049.63"
$ws.Range("D18").Value = "38.103.68"
$ws.Range("E18").Value = "  +2.46%  "
$ws.Range("D19").Value = "'6.31"
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("D20").Value = "'69.87"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("D22").Value = "'225.94"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Value = "'165.82"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").Value = "'9.24"
$ws.Range("E27").Value = "  +1.99%  "
$ws.Range("D28").Value = "'0.134"
$ws.Range("E28").Value = "  +3.78%  "
$ws.Range("D29").Value = "'18.97"
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("E30").Value = "  -1.23%  "
$ws.Range("E31").Value = "  +2.55%  "
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("E33").Value = "  +2.73%  "
$ws.Range("E34").Value = "  +8.10%  "
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("E36").Value = "  +16.15%  "
$ws.Range("E37").Value = "  -2.90%  "
$ws.Range("E38").Value = "  +2.75%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").Value = "1.518.83"
$ws.Range("E40").Value = "  +3.02%  "
$ws.Range("E41").Value = "  +4.45%  "
$ws.Range("D42").Value = "'97.60"
$ws.Range("E42").Value = "  +3.06%  "
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("E44").Value = "  +2.55%  "
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("D47").Value = "'4.06"
$ws.Range("E47").Value = "  -2.17%  "
$ws.Range("D48").Value = "'1.02"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("E49").Value = "  +1.44%  "
$ws.Range("D50").Value = "'7.02"
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("D51").Value = "2.246.99"
$ws.Range("E51").Value = "  +1.62%  "
